# Implement data parsing logic (resolves issue #62).
# Appends two new daily log rows (2025-06-11 and 2025-06-12) to each of the
# four device sheets, extending sheetData from row 32 to row 34.

$wb = $excel.ActiveWorkbook

function SetLogRow($ws, $row, $dateVal, $b, $c, $d, $e, $f, $g, $h, $i) {
    $ws.Range("A$row").Value = $dateVal
    $ws.Range("A$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e
    $ws.Range("F$row").Value = $f
    $ws.Range("G$row").Value = [double]$g
    $ws.Range("H$row").Value = $h
    $ws.Range("I$row").Value = $i
}

$newDate1 = 45819.43640046296
$newDate2 = 45820.43376157407

# --- Sheet 1: DE_LFT_#1 ---
$ws1 = $wb.Worksheets.Item(1)
SetLogRow $ws1 33 $newDate1 "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x01,0x74" "0x14" 380 "7.598631275147109e+23" 372 14
SetLogRow $ws1 34 $newDate2 "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x01,0x74" "0x14" 380 "7.598631275147109e+23" 372 14

# --- Sheet 2: DE_LFT_#2 ---
$ws2 = $wb.Worksheets.Item(2)
SetLogRow $ws2 33 $newDate1 "0x01,0x7c" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x01,0x74" "0xe" 380 "5.68432987514711e+23" 372 14
SetLogRow $ws2 34 $newDate2 "0x01,0x7c" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x01,0x74" "0xe" 380 "5.68432987514711e+23" 372 14

# --- Sheet 3: DE_PLT_#1 ---
$ws3 = $wb.Worksheets.Item(3)
SetLogRow $ws3 33 $newDate1 "0x00,0x82" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x82" "0x7" 130 "5.68631262647114e+23" 129 7
SetLogRow $ws3 34 $newDate2 "0x00,0x82" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x80" "0x7" 130 "5.68631262647114e+23" 128 7

# --- Sheet 4: DE_PLT_#2 ---
$ws4 = $wb.Worksheets.Item(4)
SetLogRow $ws4 33 $newDate1 "0x00,0x82" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x81" "0x3" 130 "9.85046333984776e+23" 129 3
SetLogRow $ws4 34 $newDate2 "0x00,0x82" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x81" "0x3" 130 "9.85046333984776e+23" 129 3
